$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph index whose text matches a given substring.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($pattern) {
  $count = $d.Paragraphs.Count
  for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match $pattern) {
      return $i
    }
  }
  return -1
}

# ---------------------------------------------------------------------------
# 1) "Note that on our Raspbian system, ..." paragraph becomes the
#    "A caveat here for Rpi users, ..." paragraph, split across 4 runs
#    that all carry the "Source Text" character style.
# ---------------------------------------------------------------------------
$idx1 = Find-ParagraphIndex("Note that on our Raspbian")
$p1 = $d.Paragraphs.Item($idx1)
$r1 = $p1.Range

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="TextBody"/>
    <w:rPr/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:color w:val="auto"/>
    </w:rPr>
    <w:t>A caveat here for Rpi users,</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:color w:val="auto"/>
    </w:rPr>
    <w:t xml:space="preserve"> on our Raspbian system, we also had to use the following command to </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:color w:val="auto"/>
    </w:rPr>
    <w:t>edit the raspi-config file, and have the boot proceed by waiting for the network to come up first</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:color w:val="auto"/>
    </w:rPr>
    <w:t>-</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r1.InsertXML($xml1)

$p1b = $d.Paragraphs.Item($idx1)
$rng1 = $p1b.Range
$textRng1 = $d.Range($rng1.Start, $rng1.End - 1)
$textRng1.Style = "Source Text"

# ---------------------------------------------------------------------------
# 2) "$ sudo mount /mnt/rasp" paragraph becomes "$ sudo raspi-config".
#    The leading "$ " run is kept, the rest collapses into "sudo " +
#    "raspi-config" (both bold), still using the "Source Text" style.
# ---------------------------------------------------------------------------
$idx2 = Find-ParagraphIndex("sudo mount /mnt/rasp")
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="TextBody"/>
    <w:rPr/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:color w:val="auto"/>
    </w:rPr>
    <w:t xml:space="preserve">$ </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="auto"/>
    </w:rPr>
    <w:t xml:space="preserve">sudo </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="auto"/>
    </w:rPr>
    <w:t>raspi-config</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r2.InsertXML($xml2)

$p2b = $d.Paragraphs.Item($idx2)
$rng2 = $p2b.Range
$textRng2 = $d.Range($rng2.Start, $rng2.End - 1)
$textRng2.Style = "Source Text"

# ---------------------------------------------------------------------------
# 3) Append a new "References:" section at the end of the document with
#    two reference entries: a plain hyperlink and a field-code hyperlink.
# ---------------------------------------------------------------------------
$docEnd = $d.Content.End
$rEnd = $d.Range($docEnd - 1, $docEnd - 1)

$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="TextBody"/>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>References:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="TextBody"/>
    <w:rPr/>
  </w:pPr>
  <w:hyperlink r:id="rIdNarrowEscape">
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
      </w:rPr>
      <w:t>Narrow Escape - Ubuntu 16.04: Share ZFS storage via NFS/SMB</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="TextBody"/>
    <w:rPr/>
  </w:pPr>
  <w:r>
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
    </w:rPr>
    <w:instrText xml:space="preserve"> HYPERLINK "https://docs.oracle.com/cd/E23824_01/html/821-1448/gayne.html" \l "NewZFSSharingSyntaxhttps://docs.oracle.com/cd/E23824_01/html/821-1448/gayne.html%23NewZFSSharingSyntax"</w:instrText>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
    </w:rPr>
    <w:fldChar w:fldCharType="separate"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
    </w:rPr>
    <w:t>Sharing and Unsharing ZFS File Systems</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="New times roman" w:hAnsi="New times roman"/>
    </w:rPr>
    <w:fldChar w:fldCharType="end"/>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdNarrowEscape" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://narrowescape.org/infotech/howtos/share-zfs-nfs-smb.html" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rEnd.InsertXML($xml3)

# Fix up paragraph spacing on the 3rd new paragraph (the field hyperlink one)
# so it carries "before=0 after=140" direct spacing, matching the target.
$newCount = $d.Paragraphs.Count
$fieldPara = $d.Paragraphs.Item($newCount)
$fieldPara.Format.SpaceBefore = 0
$fieldPara.Format.SpaceAfter = 7

# Apply the "Internet Link" character style plus the document's body font to
# both hyperlink runs (the plain hyperlink and all the field-hyperlink runs).
$refPara = $d.Paragraphs.Item($newCount - 1)
$refRng = $refPara.Range
$refTextRng = $d.Range($refRng.Start, $refRng.End - 1)
$refTextRng.Style = "Internet Link"
$refTextRng.Font.Name = "New times roman"

$fieldRng = $fieldPara.Range
$fieldTextRng = $d.Range($fieldRng.Start, $fieldRng.End - 1)
$fieldTextRng.Style = "Internet Link"
$fieldTextRng.Font.Name = "New times roman"

# ---------------------------------------------------------------------------
# 4) Remove the now-unused "ListLabel 9" character style.
# ---------------------------------------------------------------------------
$unusedStyle = $d.Styles("ListLabel 9")
$unusedStyle.Delete()

Write-Output "All edits applied."
